# Update cryptocurrency price and volume data (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '88.263.22'
$ws.Range("E2").Value = '  -0.63%  '
$ws.Range("D3").Value = '3.119.73'
$ws.Range("E3").Value = '  -1.16%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.58'
$ws.Range("E5").Value = '  +1.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '636.31'
$ws.Range("E6").Value = '  +3.99%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.388'
$ws.Range("E7").Value = '  +0.75%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.824'
$ws.Range("E8").Value = '  +20.36%  '
$ws.Range("E9").Value = '  +0.16%  '
$ws.Range("D10").Value = '3.116.08'
$ws.Range("E10").Value = '  -1.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.583'
$ws.Range("E11").Value = '  +2.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.179'
$ws.Range("E12").Value = '  +1.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000249'
$ws.Range("E13").Value = '  -1.08%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.39'
$ws.Range("E14").Value = '  +3.51%  '
$ws.Range("D15").Value = '88.038.53'
$ws.Range("E15").Value = '  -0.80%  '
$ws.Range("D16").Value = '3.691.95'
$ws.Range("E16").Value = '  -1.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '32.17'
$ws.Range("E17").Value = '  -0.73%  '
$ws.Range("D18").Value = '3.121.90'
$ws.Range("E18").Value = '  -0.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.40'
$ws.Range("E19").Value = '  +4.69%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000221'
$ws.Range("E20").Value = '  +17.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.33'
$ws.Range("E21").Value = '  +0.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '425.87'
$ws.Range("E22").Value = '  -1.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.49'
$ws.Range("E23").Value = '  -0.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.93'
$ws.Range("E24").Value = '  -1.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.64'
$ws.Range("E25").Value = '  +11.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '83.91'
$ws.Range("E26").Value = '  +12.38%  '
$ws.Range("E27").Value = '  -0.77%  '
$ws.Range("D28").Value = '3.279.05'
$ws.Range("E28").Value = '  -1.38%  '
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("E30").Value = '  -0.19%  '
$ws.Range("E31").Value = '  -5.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.93'
$ws.Range("E32").Value = '  -3.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.20'
$ws.Range("E33").Value = '  -1.87%  '
$ws.Range("E34").Value = '  +18.11%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '504.74'
$ws.Range("E35").Value = '  -4.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.86'
$ws.Range("E36").Value = '  -1.30%  '
$ws.Range("E37").Value = '  +1.84%  '
$ws.Range("E38").Value = '  -0.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.65'
$ws.Range("E39").Value = '  +4.04%  '
$ws.Range("E40").Value = '  -0.20%  '
$ws.Range("E41").Value = '  +0.55%  '
$ws.Range("E42").Value = '  -0.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.368'
$ws.Range("E43").Value = '  -0.27%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.85'
$ws.Range("E44").Value = '  -2.30%  '
$ws.Range("E45").Value = '  +11.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '146.17'
$ws.Range("E46").Value = '  -2.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '44.05'
$ws.Range("E47").Value = '  +0.88%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0676'
$ws.Range("E48").Value = '  +15.84%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '162.70'
$ws.Range("E49").Value = '  -4.60%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.718'
$ws.Range("E50").Value = '  +3.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.20'
$ws.Range("E51").Value = '  -2.05%  '
